$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "W fan_evap"
$ws.Range("G2").Value = 0.7789
$ws.Range("G3").Value = 0.7648
$ws.Range("G4").Value = 0.7801
$ws.Range("G5").Value = 0.7786
$ws.Range("G6").Value = 0.7602
$ws.Range("G7").Value = 0.7639
$ws.Range("G8").Value = 0.7696
$ws.Range("G9").Value = 0.7547

$ws.Range("G2").Select()
